$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $scratch, $text) {
    # Forces a genuinely text-typed cell value (preserving things like
    # leading zeros / decimal-looking strings) without Excel's automatic
    # "looks like a number" coercion, and without creating a new cell
    # style (unlike NumberFormat="@" or a leading quote-prefix, both of
    # which register a brand-new style record).
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues
}

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (total) sheet,
#    matching the formatting of the existing per-quarter fund-holding sheets
#    (e.g. "2021-Q4").
# ---------------------------------------------------------------------------
$totalWsAnchor = $wb.Worksheets.Item("总计")       # anchor used only to position the new sheet
$templateWs    = $wb.Worksheets.Item("2021-Q4")    # formatting template for a per-fund sheet

$ws = $wb.Worksheets.Add($totalWsAnchor)
$ws.Name = "2022-Q1"
$scratch = $ws.Range("Z1")

# Header row: copy formatting (bold / centered / bordered) from the template sheet
$templateWs.Range("B1:H1").Copy($ws.Range("B1:H1"))

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Fund holdings data rows (B/D/E/F/G are text values in the source data).
$rows = @(
    @{ idx=0; code="001150"; name="融通互联网传媒灵活配置混合"; scale="9.18"; pos="92.66"; ratio="2.94"; value="0.2699"; rank=4 },
    @{ idx=1; code="012200"; name="新华鑫科技3个月滚动持有灵活配置混合型证券投资基金A"; scale="2.04"; pos="77.02"; ratio="3.36"; value="0.0685"; rank=7 },
    @{ idx=2; code="200001"; name="长城久恒灵活配置混合"; scale="0.90"; pos="94.32"; ratio="3.01"; value="0.0271"; rank=2 },
    @{ idx=3; code="001209"; name="前海开源一带一路主题精选灵活配置混合A"; scale="0.92"; pos="82.84"; ratio="2.89"; value="0.0266"; rank=9 },
    @{ idx=4; code="001723"; name="华商新动力灵活配置混合"; scale="0.75"; pos="55.57"; ratio="2.41"; value="0.0181"; rank=7 },
    @{ idx=5; code="012201"; name="新华鑫科技3个月滚动持有灵活配置混合型证券投资基金C"; scale="0.52"; pos="77.02"; ratio="3.36"; value="0.0175"; rank=7 },
    @{ idx=6; code="004930"; name="华润元大价值优选混合A"; scale="0.32"; pos="65.19"; ratio="3.36"; value="0.0108"; rank=8 },
    @{ idx=7; code="004931"; name="华润元大价值优选混合C"; scale="0.18"; pos="65.19"; ratio="3.36"; value="0.0060"; rank=8 },
    @{ idx=8; code="620004"; name="金元顺安价值增长混合"; scale="0.18"; pos="87.28"; ratio="2.69"; value="0.0048"; rank=10 },
    @{ idx=9; code="002080"; name="前海开源一带一路主题精选灵活配置混合C"; scale="0.10"; pos="82.84"; ratio="2.89"; value="0.0029"; rank=9 }
)

$r = 2
foreach ($row in $rows) {
    # column-A index cell: copy style from the template's column-A cell, then set its numeric value
    $templateWs.Range("A2").Copy($ws.Cells.Item($r, 1))
    $ws.Cells.Item($r, 1).Value = $row.idx

    Set-TextValue $ws.Cells.Item($r, 2) $scratch $row.code
    $ws.Cells.Item($r, 3).Value = $row.name
    Set-TextValue $ws.Cells.Item($r, 4) $scratch $row.scale
    Set-TextValue $ws.Cells.Item($r, 5) $scratch $row.pos
    Set-TextValue $ws.Cells.Item($r, 6) $scratch $row.ratio
    Set-TextValue $ws.Cells.Item($r, 7) $scratch $row.value
    $ws.Cells.Item($r, 8).Value = $row.rank

    $r = $r + 1
}

$scratch.ClearContents()
$ws.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. Update the "总计" (total) sheet: insert a new first data row for
#    "2022-Q1" (10 funds, 0.45 亿元) and shift the existing rows down.
#    The worksheet reference is re-fetched by name because the previously
#    held reference gets rebound to whatever now occupies its old position
#    after a new sheet is inserted.
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

# Shift existing data rows 2..6 down to 3..7 (copy whole rows, bottom-up so
# sources aren't clobbered before they are read). Range.Copy preserves both
# values and the per-cell formatting (e.g. column A's bold/bordered style).
for ($i = 6; $i -ge 2; $i--) {
    $src = $totalWs.Range("A" + $i + ":D" + $i)
    $dst = $totalWs.Range("A" + ($i + 1) + ":D" + ($i + 1))
    $src.Copy($dst)
}

# Write the new first data row (reuses row 2's pre-existing formatting).
$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 10
$totalWs.Range("D2").Value = 0.45

# Renumber the column-A running index (0,1,2,...) for every data row now
# that the new row has shifted everything down by one.
for ($i = 3; $i -le 7; $i++) {
    $totalWs.Cells.Item($i, 1).Value = $i - 2
}

$totalWs.Range("A1").Select()
